$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$cell = $t.Cell(1,1)
$cell.Range.Text = "36 x 83" + $nl + "  8    3" + $nl + "  ----" + $nl + "3|    |" + $nl + "6|    |"

$cell = $t.Cell(1,2)
$cell.Range.Text = "41 x 59" + $nl + "  5    9" + $nl + "  ----" + $nl + "4|    |" + $nl + "1|    |"

$cell = $t.Cell(1,3)
$cell.Range.Text = "52 x 54" + $nl + "  5    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"

$cell = $t.Cell(2,1)
$cell.Range.Text = "43 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "4|    |" + $nl + "3|    |"

$cell = $t.Cell(2,2)
$cell.Range.Text = "53 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"

$cell = $t.Cell(2,3)
$cell.Range.Text = "24 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"

$cell = $t.Cell(3,1)
$cell.Range.Text = "34 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "4|    |"

$cell = $t.Cell(3,2)
$cell.Range.Text = "57 x 33" + $nl + "  3    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "7|    |"

$cell = $t.Cell(3,3)
$cell.Range.Text = "51 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "5|    |" + $nl + "1|    |"

$cell = $t.Cell(4,1)
$cell.Range.Text = "94 x 44" + $nl + "  4    4" + $nl + "  ----" + $nl + "9|    |" + $nl + "4|    |"

$cell = $t.Cell(4,2)
$cell.Range.Text = "62 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "2|    |"

$cell = $t.Cell(4,3)
$cell.Range.Text = "73 x 95" + $nl + "  9    5" + $nl + "  ----" + $nl + "7|    |" + $nl + "3|    |"

$cell = $t.Cell(5,1)
$cell.Range.Text = "80 x 16" + $nl + "  1    6" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"

$cell = $t.Cell(5,2)
$cell.Range.Text = "80 x 62" + $nl + "  6    2" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"

$cell = $t.Cell(5,3)
$cell.Range.Text = "68 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "8|    |"
